# "Prepare for new ppt"
#
# The canonical-OOXML diff for this commit shows a new, essentially
# empty notes page being attached to the first slide (ppt/notesSlides/
# notesSlide1.xml, wired up via ppt/slides/_rels/slide1.xml.rels,
# [Content_Types].xml, etc.). Re-create that with the PowerPoint
# object model: touch Slide(1).NotesPage so the notes slide part is
# materialized and persisted with the package.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Accessing .NotesPage alone only hands back an in-memory preview;
# the notes slide part is only minted (and linked into slide1's
# relationships / [Content_Types].xml / notesSlide1.xml) once the
# notes body placeholder is actually materialized, so do that
# explicitly via AddPlaceholder(2) == ppPlaceholderBody.
$notesPage = $s.NotesPage
$notesBody = $notesPage.Shapes.AddPlaceholder(2)
